$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..16 : column B (IdxPar), column C (IdxSG), column D (DesSG)
$rows = @(
    @(2,  2050000000, 305020100, "Divers"),
    @(3,  2050000000, 305020200, "Autorités d'Eglise"),
    @(4,  2050000000, 305020300, "KT/Jeunesse"),
    @(5,  2050000000, 305020400, "Bénévoles, Responsables"),
    @(6,  2050000000, 305020500, "Autorités, relations publiques"),
    @(7,  2050000000, 305020600, "Finances"),
    @(8,  2050000000, 305020700, "Information/Journal"),
    @(9,  2050000000, 305020800, "Ministres"),
    @(10, 2050000000, 305020900, "ACTIVITES NON CLASSEES"),
    @(11, 2050000000, 305021000, "Formation adultes"),
    @(12, 2050000000, 305021100, "Enfance"),
    @(13, 2050000000, 305021200, "REGISTRES PAROISSIAUX"),
    @(14, 2050000000, 305021300, "Comité visites paroissiales"),
    @(15, 2050000000, 305021400, "Mission, Solidarité"),
    @(16, 2050000000, 305021500, "Couples mixtes")
)

foreach ($r in $rows) {
    $rowIndex = $r[0]
    $bVal = $r[1]
    $cVal = $r[2]
    $dVal = $r[3]

    $bCell = $ws.Cells.Item($rowIndex, 2)
    $bCell.Value = $bVal
    # Make column B share the same cell style as the other data rows
    # (numeric / vertical-top alignment, no special text number format).
    $bCell.VerticalAlignment = -4160

    $ws.Cells.Item($rowIndex, 3).Value = $cVal
    $ws.Cells.Item($rowIndex, 4).Value = $dVal
}

# Update the active selection to match the new cursor position.
$ws.Range("H8").Select()
